$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value2 = '57.495.35'
$ws.Range('E2').Value2 = '  -4.52%  '

$ws.Range('D3').Value2 = '3.114.51'
$ws.Range('E3').Value2 = '  -5.87%  '

$ws.Range('E4').Value2 = '  -0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '520.15'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value2 = '  -6.78%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value2 = '134.14'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value2 = '  -5.34%  '

$ws.Range('E7').Value2 = '  -0.12%  '

$ws.Range('D8').Value2 = '3.111.84'
$ws.Range('E8').Value2 = '  -6.03%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '0.442'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value2 = '  -6.63%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '7.18'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value2 = '  -8.60%  '

$ws.Range('E11').Value2 = '  -8.07%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value2 = '0.382'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value2 = '  -6.10%  '

$ws.Range('D13').Value2 = '3.649.66'
$ws.Range('E13').Value2 = '  -5.92%  '

$ws.Range('E14').Value2 = '  -2.31%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '25.30'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value2 = '  -5.56%  '

$ws.Range('D16').Value2 = '3.120.39'
$ws.Range('E16').Value2 = '  -5.60%  '

$ws.Range('D17').Value2 = '57.442.97'
$ws.Range('E17').Value2 = '  -4.64%  '

$ws.Range('E18').Value2 = '  -8.72%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '5.75'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value2 = '  -6.80%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '12.91'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value2 = '  -10.32%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value2 = '7.94'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value2 = '  -7.97%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '340.85'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value2 = '  -9.09%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value2 = '0.999'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value2 = '  -0.14%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value2 = '68.11'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value2 = '  -8.24%  '

$ws.Range('E25').Value2 = '  -7.52%  '

$ws.Range('D26').Value2 = '3.247.78'
$ws.Range('E26').Value2 = '  -5.74%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value2 = '0.165'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value2 = '  -4.05%  '

$ws.Range('E28').Value2 = '  -0.18%  '

$ws.Range('D29').Value2 = '0.0₃0935'
$ws.Range('E29').Value2 = '  -9.06%  '

$ws.Range('E30').Value2 = '  -0.22%  '

$ws.Range('E31').Value2 = '  -6.73%  '

$ws.Range('E32').Value2 = '  -8.58%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value2 = '6.88'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value2 = '  -10.04%  '

$ws.Range('B34').Value2 = 'EthereumClassic'
$ws.Range('C34').Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value2 = '21.35'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value2 = '  -5.26%  '

$ws.Range('B35').Value2 = 'Fetch.AI'
$ws.Range('C35').Value2 = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value2 = '1.22'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value2 = '  -3.61%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value2 = '157.89'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value2 = '  -4.79%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value2 = '4.75'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value2 = '  -7.73%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value2 = '6.13'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value2 = '  -9.15%  '

$ws.Range('E39').Value2 = '  -10.87%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '25.08'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value2 = '  -6.29%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '0.0686'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value2 = '  -7.04%  '

$ws.Range('D42').Value2 = '3.141.63'
$ws.Range('E42').Value2 = '  -5.93%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '40.26'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value2 = '  -4.04%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '0.681'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value2 = '  -9.60%  '

$ws.Range('E45').Value2 = '  -4.72%  '

$ws.Range('E46').Value2 = '  -7.42%  '

$ws.Range('E47').Value2 = '  +0.04%  '

$ws.Range('E48').Value2 = '  -9.56%  '

$ws.Range('D49').Value2 = '2.247.33'
$ws.Range('E49').Value2 = '  -4.82%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '6.13'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value2 = '  -5.96%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value2 = '19.85'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value2 = '  -6.59%  '
